# Insert a new data row at row 318 (weekly price update), pushing the
# existing rows 318-347 down to 319-348. This mirrors the author's edit:
# a brand-new "Repollo" price observation for "Feria Lagunitas de Puerto
# Montt" is inserted into the middle of the table, and every subsequent
# row shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 318 - this shifts rows
# 318..347 down to 319..348 (and Excel auto-extends the sheet dimension).
$ws.Rows.Item(318).Insert()

# Populate the newly inserted row 318 with the new record's values.
$ws.Range("A318").Value = 4
$ws.Range("B318").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C318").Value = "Los Lagos"
$ws.Range("D318").Value = 44578
$ws.Range("E318").Value = 10
$ws.Range("F318").Value = 100112006
$ws.Range("G318").Value = "Repollo"
$ws.Range("H318").Value = "Crespo record"
$ws.Range("I318").Value = "Segunda"
$ws.Range("J318").Value = 500
$ws.Range("K318").Value = 1000
$ws.Range("L318").Value = 1000
$ws.Range("M318").Value = 1000
$ws.Range("N318").Value = "$/unidad"
$ws.Range("O318").Value = "Región del Maule"
$ws.Range("P318").Value = 1000
$ws.Range("Q318").Value = 1
$ws.Range("R318").Value = "Hortaliza"
